$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(" Abu Dhabi", " October 30 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Rajasthan Royals", "Kings XI Punjab", "Sanju Samson †", "48", "25", "4", "3", "192.00"),
    @(" Abu Dhabi", " October 25 2020", "Royals won by 8 wickets (with 10 balls remaining)", "Rajasthan Royals", "Mumbai Indians", "Sanju Samson †", "54", "31", "4", "3", "174.19"),
    @(" Dubai (DSC)", " October 22 2020", "Sunrisers won by 8 wickets (with 11 balls remaining)", "Rajasthan Royals", "Sunrisers Hyderabad", "Sanju Samson †", "36", "26", "3", "1", "138.46"),
    @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Rajasthan Royals", "Chennai Super Kings", "Sanju Samson †", "74", "32", "1", "9", "231.25"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Rajasthan Royals", "Kolkata Knight Riders", "Sanju Samson †", "1", "4", "0", "0", "25.00"),
    @(" Abu Dhabi", " October 19 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Rajasthan Royals", "Chennai Super Kings", "Sanju Samson †", "0", "3", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 17 2020", "RCB won by 7 wickets (with 2 balls remaining)", "Rajasthan Royals", "Royal Challengers Bangalore", "Sanju Samson †", "9", "6", "0", "1", "150.00")
)

$startRow = 9
$endRow = $startRow + $rows.Count - 1

# Columns G:K hold numeric-looking values that must stay stored as plain
# text (matching the rest of the sheet), so force text format before
# writing so COM does not auto-convert them into real numbers.
$ws.Range("G$startRow`:K$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
